# Lecture 02 - Image Filtering.pptx
# Update the deck from "Lecture 01" to "Lecture 02":
#  - Title slide: "Image Filtering" -> "Lecture 02 - Image Filtering"
#  - Every recurring footer shape: "Lecture 01" -> "Lecture 02"

$p = $ppt.ActivePresentation

# --- Slide 1: title shape gets a "Lecture 02 - " prefix -------------------
$titleSlide = $p.Slides.Item(1)
for ($i = 1; $i -le $titleSlide.Shapes.Count; $i++) {
    $sh = $titleSlide.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "Image Filtering") {
        $tr = $sh.TextFrame.TextRange
        $null = $tr.InsertBefore(" 02 - ")
        $null = $tr.InsertBefore("Lecture")
        break
    }
}

# --- All slides: footer shape "Lecture 01" -> "Lecture 02" ----------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "Lecture 01") {
                $found = $tr.Find(" 01", 0)
                if ($found -ne $null) {
                    $sub = $tr.Characters($found.Start, $found.Length)
                    $sub.Text = " 02"
                }
            }
        }
    }
}
